$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the CGM test case id (row 2) to a newly generated id, as part of the
# "update genetics, cgm, eeg" refresh of generated case ids.
$ws.Range("A2").Value = "CA-756V081T"
